$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("List1")

# Improved startup var calculation: update the rated-speed (G2), the startup
# inertia/time constant (J2), and have G5 reference J2 instead of a hard-coded
# literal so downstream formulas (H5, I5, A9:D9, A15:C15) stay in sync.
$ws1.Range("G2").Value = 60000
$ws1.Range("J2").Value = 0.00005
$ws1.Range("G5").Formula = "=F5*J2"

# Selection / active sheet changes: List1 becomes the active (tabSelected) sheet,
# with the selection moved from F6 to G6
$ws1.Activate()
$ws1.Range("G6").Select()
